# Rerun some analyses with more data: append the three newly-added
# datasets (MOT, SOCON, LISS) as additional rows below the existing
# dataset/wave/year/color table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New records, in the same column order as the existing table:
# dataset, wave, year, color
$newRows = @(
    @("MOT",   1, 2019, "#E02023"),
    @("MOT",   2, 2021, "#E02023"),
    @("SOCON", 1, 2020, "#1676E3"),
    @("SOCON", 2, 2021, "#1676E3"),
    @("SOCON", 3, 2022, "#1676E3"),
    @("LISS",  1, 2019, "#0AC760"),
    @("LISS",  2, 2020, "#0AC760"),
    @("LISS",  3, 2020, "#0AC760"),
    @("LISS",  4, 2021, "#0AC760")
)

$startRow = 26

# Fill column-by-column (all of A, then all of B, then C, then D) so
# that the new shared-string entries are interned in dataset-name order
# first (MOT, SOCON, LISS) and then color order (#E02023, #1676E3,
# #0AC760), matching how the workbook was actually authored.
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Range("A$r").Value = $newRows[$i][0]
}
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Range("B$r").Value = $newRows[$i][1]
}
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Range("C$r").Value = $newRows[$i][2]
}
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Range("D$r").Value = $newRows[$i][3]
}

# Leave the cursor where the author left it when they saved the file.
$ws.Range("E32").Select() | Out-Null
